$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "E"="3"; "G"="2.475141333333333"; "H"="7.425424"; "I"="0.004126344150637133"; "J"="0.004126344150637133"; "K"="3"; "M"="3.433022333333334"; "N"="10.299067"; "O"="0.2851620693776887"; "P"="0.2851620693776887"; "Q"="8.497215475489778"; "R"="76.47493927940801"; "S"="0.001176676836960206"; "T"="0.001176676836960206" }
    3 = @{ "E"="3"; "G"="2.475141333333333"; "H"="7.425424"; "I"="0.004126344150637133"; "J"="0.004126344150637133"; "K"="3"; "M"="8.017154666666668"; "N"="24.051464"; "O"="0.6659404435181344"; "P"="0.6659404435181344"; "Q"="19.84359089119289"; "R"="178.592318020736"; "S"="0.002747899453783752"; "T"="0.002747899453783752" }
    4 = @{ "E"="3"; "G"="2.475141333333333"; "H"="7.425424"; "I"="0.004126344150637133"; "J"="0.004126344150637133"; "K"="3"; "M"="0.5886693333333334"; "N"="1.766008"; "O"="0.04889748710417684"; "P"="0.04889748710417684"; "Q"="1.457039798599111"; "R"="13.113358187392"; "S"="0.0002017678598931747"; "T"="0.0002017678598931747" }
    5 = @{ "E"="3"; "G"="577.349335"; "H"="1732.048005"; "I"="0.9625074816003053"; "J"="0.9625074816003053"; "K"="3"; "M"="3.433022333333334"; "N"="10.299067"; "O"="0.2851620693776887"; "P"="0.2851620693776887"; "Q"="1982.053161190149"; "R"="17838.47845071134"; "S"="0.2744706252446507"; "T"="0.2744706252446507" }
    6 = @{ "E"="3"; "G"="577.349335"; "H"="1732.048005"; "I"="0.9625074816003053"; "J"="0.9625074816003053"; "K"="3"; "M"="8.017154666666668"; "N"="24.051464"; "O"="0.6659404435181344"; "P"="0.6659404435181344"; "Q"="4628.698915392148"; "R"="41658.29023852933"; "S"="0.6409726591864299"; "T"="0.6409726591864299" }
    7 = @{ "E"="3"; "G"="577.349335"; "H"="1732.048005"; "I"="0.9625074816003053"; "J"="0.9625074816003053"; "K"="3"; "M"="0.5886693333333334"; "N"="1.766008"; "O"="0.04889748710417684"; "P"="0.04889748710417684"; "Q"="339.8678481348933"; "R"="3058.81063321404"; "S"="0.04706419716922465"; "T"="0.04706419716922465" }
    8 = @{ "E"="3"; "G"="20.014326"; "H"="60.04297800000001"; "I"="0.03336617424905757"; "J"="0.03336617424905757"; "K"="3"; "M"="3.433022333333334"; "N"="10.299067"; "O"="0.2851620693776887"; "P"="0.2851620693776887"; "Q"="68.709628144614"; "R"="618.3866533015261"; "S"="0.009514767296077805"; "T"="0.009514767296077807" }
    9 = @{ "E"="3"; "G"="20.014326"; "H"="60.04297800000001"; "I"="0.03336617424905757"; "J"="0.03336617424905757"; "K"="3"; "M"="8.017154666666668"; "N"="24.051464"; "O"="0.6659404435181344"; "P"="0.6659404435181344"; "Q"="160.457947091088"; "R"="1444.121523819792"; "S"="0.02221988487792076"; "T"="0.02221988487792076" }
    10 = @{ "E"="3"; "G"="20.014326"; "H"="60.04297800000001"; "I"="0.03336617424905757"; "J"="0.03336617424905757"; "K"="3"; "M"="0.5886693333333334"; "N"="1.766008"; "O"="0.04889748710417684"; "P"="0.04889748710417684"; "Q"="11.781819943536"; "R"="106.036379491824"; "S"="0.00163152207505901"; "T"="0.00163152207505901" }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = [double]$cols[$col]
    }
}
